# Student & Hall Modified
# Adds a "Sex" column (D) to the student details table:
#   D4 = "Sex" header (bold, size 16)
#   D5:D8 = "M" for every student row
# Also bumps row 4's height to fit the bigger header font and moves the
# saved cell selection, matching the authored workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell for the "Sex" column.
$ws.Range("D4").Value = "Sex"
$ws.Range("D4").Font.Bold = $true
$ws.Range("D4").Font.Size = 16

# Row 4 grows a bit to accommodate the larger header font.
$ws.Rows.Item(4).RowHeight = 21

# Every student is marked "M" in the new column.
$ws.Range("D5").Value = "M"
$ws.Range("D6").Value = "M"
$ws.Range("D7").Value = "M"
$ws.Range("D8").Value = "M"

# Match the saved selection left in the workbook.
$ws.Range("E16").Select()
